$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts existing rows 21-28 down to 22-29),
# then populate it with a new weekly price entry (same record as the old
# row 21, just dated one week later).
$ws.Rows.Item(21).Insert()

$ws.Cells.Item(21, 1).Value()  = 4
$ws.Cells.Item(21, 2).Value()  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value()  = "Los Lagos"
$ws.Cells.Item(21, 4).Value()  = 44754
$ws.Cells.Item(21, 5).Value()  = 10
$ws.Cells.Item(21, 6).Value()  = 100112012
$ws.Cells.Item(21, 7).Value()  = "Espinaca"
$ws.Cells.Item(21, 8).Value()  = "Sin especificar"
$ws.Cells.Item(21, 9).Value()  = "Primera"
$ws.Cells.Item(21, 10).Value() = 30
$ws.Cells.Item(21, 11).Value() = 13000
$ws.Cells.Item(21, 12).Value() = 13000
$ws.Cells.Item(21, 13).Value() = 13000
$ws.Cells.Item(21, 14).Value() = "$/cuna 10 kilos"
$ws.Cells.Item(21, 15).Value() = "Región Metropolitana"
$ws.Cells.Item(21, 16).Value() = 1300
$ws.Cells.Item(21, 17).Value() = 10
$ws.Cells.Item(21, 18).Value() = "Hortaliza"
